# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to match the newly generated site output (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 136
$ws1.Range("F3").Value  = 2170
$ws1.Range("F4").Value  = 46
$ws1.Range("F5").Value  = 11452
$ws1.Range("F7").Value  = 319
$ws1.Range("F9").Value  = 11393
$ws1.Range("F10").Value = 463
$ws1.Range("F11").Value = 1157
$ws1.Range("F13").Value = 1748
$ws1.Range("F14").Value = 5671
$ws1.Range("F15").Value = 108
$ws1.Range("F16").Value = 3485
$ws1.Range("F18").Value = 13

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 136
$ws4.Range("F3").Value  = 2170
$ws4.Range("F5").Value  = 46
$ws4.Range("F7").Value  = 11452
$ws4.Range("F9").Value  = 319
$ws4.Range("F11").Value = 11393
$ws4.Range("F12").Value = 463
$ws4.Range("F13").Value = 1157
$ws4.Range("F15").Value = 1748
$ws4.Range("F17").Value = 5671
$ws4.Range("F18").Value = 108
$ws4.Range("F19").Value = 3485
$ws4.Range("F21").Value = 13
